# Apply the "12th Commit with ExcelCorrection" edits to TestData.xlsx
$wb = $excel.ActiveWorkbook

$wsSignIn = $wb.Worksheets.Item("SignIn")
$wsCreate = $wb.Worksheets.Item("CreateAccount")

# ---------------------------------------------------------------------------
# 1. Header-row highlight style (new fill, theme 3 / tint 0.4) on row 1 of
#    both sheets.
# ---------------------------------------------------------------------------
$headerSignIn = $wsSignIn.Range("A1:D1")
$headerSignIn.Interior.ThemeColor = 3
$headerSignIn.Interior.TintAndShade = 0.39997558519241921

$headerCreate = $wsCreate.Range("A1:O1")
$headerCreate.Interior.ThemeColor = 3
$headerCreate.Interior.TintAndShade = 0.39997558519241921

# ---------------------------------------------------------------------------
# 2. Fix the typo'd e-mail address used for the Email / hyperlink columns.
#    SignIn!C2:C3 and CreateAccount!F2:F3.
# ---------------------------------------------------------------------------
$wsSignIn.Range("C2").Value = "testjaga14717@gmail.com"
$wsSignIn.Range("C3").Value = "testjaga14717@gmail.com"

$wsCreate.Range("F2").Value = "testjaga14717@gmail.com"
$wsCreate.Range("F3").Value = "testjaga14717@gmail.com"

# ---------------------------------------------------------------------------
# 3. CreateAccount sheet data corrections.
# ---------------------------------------------------------------------------
$wsCreate.Range("G2").Value = "IT Company"
$wsCreate.Range("G3").Value = "IT Company"

$wsCreate.Range("J2").Value = "chennai"
$wsCreate.Range("J3").Value = "chennai"

$wsCreate.Range("K2").Value = "TamilNadu"
$wsCreate.Range("K3").Value = "TamilNadu"

$wsCreate.Range("H3").Value = "Sholinganallur"

# Zipcode: store as text "600119" in both rows (was numeric 600119 / 600120)
$wsCreate.Range("I2").ClearFormats()
$wsCreate.Range("I2").Value = "'600119"
$wsCreate.Range("I3").ClearFormats()
$wsCreate.Range("I3").Value = "'600119"

# TelPhoneNo: store as text "9876543211" in both rows (was numeric) while
# keeping the existing cell style (numFmtId 49, no quote-prefix).
$wsCreate.Range("M2").NumberFormat = "@"
$wsCreate.Range("M2").Value = "9876543211"
$wsCreate.Range("M3").NumberFormat = "@"
$wsCreate.Range("M3").Value = "9876543211"

# ---------------------------------------------------------------------------
# 4. Sheet tab selection moves from CreateAccount to SignIn.
# ---------------------------------------------------------------------------
$wsCreate.Range("D8").Select()
$wsSignIn.Range("E13").Select()
$wsSignIn.Activate()

$wb.Save()
